# GBDS DECEMBER FILES 2025 - fliqlo@GBDS
# Renames the working sheet from 12-07-2025 to 12-03-2025, fills in the
# ordered-quantity figures for PP 1000 / RH 500 / RH 1000, and updates the
# active-sheet view (scrolled + selection) to match where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename the sheet tab, then re-point the Print_Area defined name (Excel
#    would normally keep this in sync on a manual tab rename, but we set it
#    explicitly to be safe).
$ws.Name = "12-03-2025"
$ws.PageSetup.PrintArea = '$A$1:$V$97'

# 2) Fill in the order quantities that were entered for this date.
$ws.Range("M13").Value = 216
$ws.Range("M15").Value = 96
$ws.Range("M16").Value = 1512

# 3) Update the view: scroll down so row 42 is at the top, and move the
#    selection in the right-hand (unfrozen) pane from M14 to M16.
$ws.Activate()
$ws.Range("M16").Select()
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 6
